$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TG102V")
$dst = $ws.Range("C25")
$dst.Value = 44391
$src = $ws.Range("C21")
$src.Copy()
$dst.PasteSpecial(-4122)
Write-Host "done"
